$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5's B,C,D,E,F,H cells were empty placeholder cells; clear them fully so they
# no longer persist as (empty) cell entries, matching the target sheet which drops
# them and keeps only the A5/G5 values.
$ws.Range("B5:F5").ClearContents()
$ws.Range("H5").ClearContents()

# New rows 6 and 7 hold reference data for two 5-digit part numbers. The rest of
# the sheet stores its numeric-looking values as text, so force text formatting
# before writing (otherwise Excel would coerce these into numbers), then drop the
# explicit style back to Normal so the cells don't carry a lingering format.
$newRows = $ws.Range("A6:H7")
$newRows.NumberFormat = "@"

$ws.Range("A6").Value = "11111"
$ws.Range("B6").Value = "1"
$ws.Range("C6").Value = "2"
$ws.Range("D6").Value = "3"
$ws.Range("E6").Value = "4"
$ws.Range("F6").Value = "5"
$ws.Range("G6").Value = "6"
$ws.Range("H6").Value = "1"

$ws.Range("A7").Value = "11112"
$ws.Range("B7").Value = "6"
$ws.Range("C7").Value = "2"
$ws.Range("D7").Value = "4"
$ws.Range("E7").Value = "4"
$ws.Range("F7").Value = "1"
$ws.Range("G7").Value = "ventilador, chapa trasera"
$ws.Range("H7").Value = "1"

$newRows.Style = "Normal"
